$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Save" in H1. Copy the formatting already used by the other
# header cells (e.g. G1: bold, centered, bordered) onto the new header cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# "Save" flag values for each data row (1 = a save was recorded for that
# game, 0 = otherwise)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
